# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# (percentages, dollar amounts, large numbers) in specific bullet / impact
# paragraphs, matching the author's "quantitative metrics highlighting"
# commit.
#
# NOTE: this interpreter only supports positional function arguments, so
# Highlight-Metrics is declared/called without named parameters.

$d = $word.ActiveDocument

function Highlight-Metrics($ParaIndex, $Metrics) {
    $para = $d.Paragraphs.Item($ParaIndex)
    $paraRange = $para.Range

    foreach ($metric in $Metrics) {
        $findRange = $paraRange.Duplicate
        $findRange.Find.ClearFormatting()
        $findRange.Find.Text = $metric
        $findRange.Find.MatchCase = $true
        $findRange.Find.MatchWholeWord = $false
        $findRange.Find.MatchWildcards = $false
        $findRange.Find.Forward = $true
        $findRange.Find.Wrap = 0
        $findRange.Find.Execute() | Out-Null

        if ($findRange.Find.Found) {
            $findRange.Font.Bold = 1
            $findRange.Font.Color = 5258796   # RGB 2C3E50 (BGR-packed int used by this host)
        }
    }
}

# Partner - Siege Analytics: race coding errors / demographic classification bullet
Highlight-Metrics 9 @("23%", "64%")

# Partner - Siege Analytics: voter turnout prediction accuracy bullet (long version)
Highlight-Metrics 11 @("87%", "71%", "±4.2%", "±2.1%")

# Senior Analyst - Myers Research: RFP vendor bids bullet
Highlight-Metrics 31 @("1,200")

# Programmer - Lake Research Partners: Polling Consortium Database bullet
Highlight-Metrics 46 @("`$400M", "`$1B")

# KEY ACHIEVEMENTS AND IMPACT: mapping cost reduction bullet
Highlight-Metrics 63 @("73.5%", "`$4.7M")

# KEY ACHIEVEMENTS AND IMPACT: voter turnout prediction accuracy bullet (short version)
Highlight-Metrics 65 @("87%", "71%")

Write-Output "Highlighted quantitative metrics in 6 paragraphs"
